$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their exact text representation
# (values such as "97.50" or "43.805.63" must not be reinterpreted as numbers)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range('D2').Value = '43.805.63'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '2.295.79'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '97.50'
$ws.Range('E5').Value = '  +2.81%  '
$ws.Range('D6').Value = '268.77'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').Value = '0.610'
$ws.Range('E9').Value = '  -1.83%  '
$ws.Range('D10').Value = '45.67'
$ws.Range('E10').Value = '  +2.00%  '
$ws.Range('D11').Value = '0.0935'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = '7.92'
$ws.Range('E12').Value = '  -2.15%  '
$ws.Range('E13').Value = '  +1.33%  '
$ws.Range('D14').Value = '2.638.86'
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('D15').Value = '15.54'
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('D16').Value = '0.854'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').Value = '2.288.33'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').Value = '43.706.47'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('E19').Value = '  +4.18%  '
$ws.Range('D20').Value = '6.19'
$ws.Range('E20').Value = '  -1.90%  '
$ws.Range('D21').Value = '71.96'
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('E22').Value = '  +11.86%  '
$ws.Range('D23').Value = '233.14'
$ws.Range('E23').Value = '  -2.09%  '
$ws.Range('E24').Value = '  -4.62%  '
$ws.Range('D25').Value = '2.67'
$ws.Range('E25').Value = '  +6.81%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '11.27'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('E28').Value = '  +2.09%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.29'
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '39.33'
$ws.Range('E30').Value = '  +1.90%  '
$ws.Range('D31').Value = '175.04'
$ws.Range('E31').Value = '  +2.02%  '
$ws.Range('D32').Value = '21.95'
$ws.Range('E32').Value = '  -2.49%  '
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('D34').Value = '5.40'
$ws.Range('E34').Value = '  -2.36%  '
$ws.Range('D35').Value = '0.126'
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '4.45'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '0.108'
$ws.Range('E37').Value = '  -0.87%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.0352'
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('E39').Value = '  -1.73%  '
$ws.Range('E40').Value = '  +2.41%  '
$ws.Range('E41').Value = '  +0.85%  '
$ws.Range('D42').Value = '12.31'
$ws.Range('E42').Value = '  +1.86%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').Value = '64.59'
$ws.Range('E44').Value = '  +4.61%  '
$ws.Range('D45').Value = '8.80'
$ws.Range('E45').Value = '  -3.05%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '0.102'
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').Value = '5.15'
$ws.Range('E47').Value = '  -5.58%  '
$ws.Range('D48').Value = '97.67'
$ws.Range('E48').Value = '  -2.56%  '
$ws.Range('E49').Value = '  -1.17%  '
$ws.Range('D50').Value = '1.51'
$ws.Range('E50').Value = '  +12.13%  '
$ws.Range('D51').Value = '2.517.77'
$ws.Range('E51').Value = '  -0.58%  '
